$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 / B3: replace placeholder "NAN" text values with numeric 100
$ws.Range("B2").Value = 100
$ws.Range("B3").Value = 100

# B4 / B5: update formulas to divide by (12*150) instead of (12*1000)
$ws.Range("B4").Formula = "=B2/(12*150)"
$ws.Range("B5").Formula = "=B3/(12*150)"

# B6 / B7: updated computed values
$ws.Range("B6").Value = 0.096541412416107
$ws.Range("B7").Value = 19.675137572927
